$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells whose new text would otherwise be auto-parsed as a number
# by Excel (single decimal point); force Text format, assign, then restore the
# default "Normal" style so no stray number-format style lingers on the cell.

$ws.Range("D2").Value = "62.891.46"
$ws.Range("E2").Value = "  +2.11%  "

$ws.Range("D3").Value = "3.480.12"
$ws.Range("E3").Value = "  +2.52%  "

$ws.Range("E4").Value = "  +0.04%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "582.65"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "147.57"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.40%  "

$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  +2.18%  "

$ws.Range("E11").Value = "  +3.58%  "

$ws.Range("D12").Value = "4.076.79"
$ws.Range("E12").Value = "  +2.57%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "29.94"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +5.15%  "

$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").Value = "3.482.36"
$ws.Range("E15").Value = "  +2.19%  "

$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").Value = "62.928.50"
$ws.Range("E17").Value = "  +2.14%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "6.35"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +3.08%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "14.40"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +5.18%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "9.34"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +3.91%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "390.13"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.567"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "75.20"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("E25").Value = "  +2.58%  "

$ws.Range("E26").Value = "  +2.83%  "

$ws.Range("E27").Value = "  -6.80%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.68"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +5.60%  "

$ws.Range("E29").Value = "  +0.13%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "8.24"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "

$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("E33").Value = "  +0.03%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "23.86"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +2.21%  "

$ws.Range("E35").Value = "  +2.85%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.28"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.37%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "31.78"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +22.34%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "171.54"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +2.00%  "

$ws.Range("E39").Value = "  +7.24%  "

$ws.Range("E40").Value = "  +2.69%  "

$ws.Range("E41").Value = "  -0.43%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.805"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +3.08%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "42.27"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("E45").Value = "  +3.35%  "

$ws.Range("E46").Value = "  +3.95%  "

$ws.Range("D47").Value = "2.609.97"
$ws.Range("E47").Value = "  +5.46%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "23.59"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.16%  "

$ws.Range("E49").Value = "  +9.16%  "

$ws.Range("E51").Value = "  -0.01%  "
